# Drafted Sessions 8 and 9
# - Watkins, Heather Jo (row 49) marked as "Sick" for the session in column C.
# - Zoe's (row 50) continuation-row session updated from "5.1 Calarco" to
#   "9.2 Martin-Caughey".
# - Scroll/selection state updated to match where the author left the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C49").Value = "Sick"
$ws.Range("C50").Value = "9.2 Martin-Caughey"

$ws.Range("A47").Select()
